# This script applies the scheduled data-refresh update produced by the
# automated Sheets runner. For each affected worksheet/row it rewrites the
# recalculated market columns (H, I, J, K, L, M, N) with their new values.
# A few rows gain or lose the "M" (profit) cell entirely when the refreshed
# computation no longer produces/now produces a value for it.
$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 1054.6346  # H15: 1236.5571 -> 1054.6346
$ws.Cells.Item(15, 9).Value = 1054.6346  # I15: 1236.5571 -> 1054.6346
$ws.Cells.Item(15, 11).Value = 3163.9038  # K15: 3709.6713 -> 3163.9038
$ws.Cells.Item(15, 13).Value = -2994.9038  # M15: -3540.6713 -> -2994.9038
# Row 33
$ws.Cells.Item(33, 8).Value = 377.13333  # H33: 356.625 -> 377.13333
$ws.Cells.Item(33, 9).Value = 283.91666  # I33: 265.84616 -> 283.91666
$ws.Cells.Item(33, 11).Value = 283.91666  # K33: 265.84616 -> 283.91666
$ws.Cells.Item(33, 13).Value = -54.91665999999998  # M33: -36.84616 -> -54.91665999999998
# Row 62
$ws.Cells.Item(62, 8).Value = 5026.6665  # H62: 2896.7715 -> 5026.6665
$ws.Cells.Item(62, 9).Value = 5428.846  # I62: 2932.9 -> 5428.846
$ws.Cells.Item(62, 10).Value = 3981  # J62: 2680 -> 3981
$ws.Cells.Item(62, 11).Value = 5428.846  # K62: 2932.9 -> 5428.846
$ws.Cells.Item(62, 12).Value = 3981  # L62: 2680 -> 3981
$ws.Cells.Item(62, 13).Value = -4804.846  # M62: -2308.9 -> -4804.846
$ws.Cells.Item(62, 14).Value = -5229  # N62: -3928 -> -5229
# Row 65
$ws.Cells.Item(65, 8).Value = 5026.6665  # H65: 2896.7715 -> 5026.6665
$ws.Cells.Item(65, 9).Value = 5428.846  # I65: 2932.9 -> 5428.846
$ws.Cells.Item(65, 10).Value = 3981  # J65: 2680 -> 3981
$ws.Cells.Item(65, 11).Value = 27144.23  # K65: 14664.5 -> 27144.23
$ws.Cells.Item(65, 12).Value = 19905  # L65: 13400 -> 19905
$ws.Cells.Item(65, 13).Value = -24024.23  # M65: -11544.5 -> -24024.23
$ws.Cells.Item(65, 14).Value = -26145  # N65: -19640 -> -26145
# Row 80
$ws.Cells.Item(80, 8).Value = 4649.6924  # H80: 4183.793 -> 4649.6924
$ws.Cells.Item(80, 9).Value = 399.58334  # I80: 360.64285 -> 399.58334
$ws.Cells.Item(80, 10).Value = 8292.643  # J80: 7752.067 -> 8292.643
$ws.Cells.Item(80, 11).Value = 1198.75002  # K80: 1081.92855 -> 1198.75002
$ws.Cells.Item(80, 12).Value = 24877.929  # L80: 23256.201 -> 24877.929
$ws.Cells.Item(80, 13).Value = -200.7500199999999  # M80: -83.92855000000009 -> -200.7500199999999
$ws.Cells.Item(80, 14).Value = -26873.929  # N80: -25252.201 -> -26873.929
# Row 83
$ws.Cells.Item(83, 8).Value = 4649.6924  # H83: 4183.793 -> 4649.6924
$ws.Cells.Item(83, 9).Value = 399.58334  # I83: 360.64285 -> 399.58334
$ws.Cells.Item(83, 10).Value = 8292.643  # J83: 7752.067 -> 8292.643
$ws.Cells.Item(83, 11).Value = 3596.25006  # K83: 3245.78565 -> 3596.25006
$ws.Cells.Item(83, 12).Value = 74633.787  # L83: 69768.603 -> 74633.787
$ws.Cells.Item(83, 13).Value = 1395.74994  # M83: 1746.21435 -> 1395.74994
$ws.Cells.Item(83, 14).Value = -84617.787  # N83: -79752.603 -> -84617.787
# Row 86
$ws.Cells.Item(86, 8).Value = 1794.2858  # H86: 1794.3334 -> 1794.2858
$ws.Cells.Item(86, 9).Value = 1537.2  # I86: 1537.2667 -> 1537.2
$ws.Cells.Item(86, 11).Value = 1537.2  # K86: 1537.2667 -> 1537.2
$ws.Cells.Item(86, 13).Value = -414.2  # M86: -414.2666999999999 -> -414.2
# Row 89
$ws.Cells.Item(89, 8).Value = 1794.2858  # H89: 1794.3334 -> 1794.2858
$ws.Cells.Item(89, 9).Value = 1537.2  # I89: 1537.2667 -> 1537.2
$ws.Cells.Item(89, 11).Value = 7686  # K89: 7686.3335 -> 7686
$ws.Cells.Item(89, 13).Value = -2070  # M89: -2070.3335 -> -2070
# Row 98
$ws.Cells.Item(98, 8).Value = 1742.75  # H98: 2177.0435 -> 1742.75
$ws.Cells.Item(98, 9).Value = 1699.4445  # I98: 2003.2727 -> 1699.4445
$ws.Cells.Item(98, 10).Value = 2912  # J98: 6000 -> 2912
$ws.Cells.Item(98, 11).Value = 1699.4445  # K98: 2003.2727 -> 1699.4445
$ws.Cells.Item(98, 12).Value = 2912  # L98: 6000 -> 2912
$ws.Cells.Item(98, 13).Value = -201.4445000000001  # M98: -505.2727 -> -201.4445000000001
$ws.Cells.Item(98, 14).Value = -5908  # N98: -8996 -> -5908
# Row 100
$ws.Cells.Item(100, 8).Value = 1505.3636  # H100: 1376.4615 -> 1505.3636
$ws.Cells.Item(100, 9).Value = 1384.3334  # I100: 1254 -> 1384.3334
$ws.Cells.Item(100, 11).Value = 1384.3334  # K100: 1254 -> 1384.3334
$ws.Cells.Item(100, 13).Value = -843.3334  # M100: -713 -> -843.3334
# Row 122
$ws.Cells.Item(122, 8).Value = 1742.75  # H122: 2177.0435 -> 1742.75
$ws.Cells.Item(122, 9).Value = 1699.4445  # I122: 2003.2727 -> 1699.4445
$ws.Cells.Item(122, 10).Value = 2912  # J122: 6000 -> 2912
$ws.Cells.Item(122, 11).Value = 5098.333500000001  # K122: 6009.8181 -> 5098.333500000001
$ws.Cells.Item(122, 12).Value = 8736  # L122: 18000 -> 8736
$ws.Cells.Item(122, 13).Value = -2648.333500000001  # M122: -3559.8181 -> -2648.333500000001
$ws.Cells.Item(122, 14).Value = -13636  # N122: -22900 -> -13636
# Row 137
$ws.Cells.Item(137, 8).Value = 13984.88  # H137: 14530.292 -> 13984.88
$ws.Cells.Item(137, 9).Value = 17380.37  # I137: 18296.223 -> 17380.37
$ws.Cells.Item(137, 11).Value = 52141.11  # K137: 54888.66900000001 -> 52141.11
$ws.Cells.Item(137, 13).Value = -49591.11  # M137: -52338.66900000001 -> -49591.11

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 33721.695  # H32: 34763.406 -> 33721.695
$ws.Cells.Item(32, 9).Value = 33721.695  # I32: 34763.406 -> 33721.695
$ws.Cells.Item(32, 11).Value = 33721.695  # K32: 34763.406 -> 33721.695
$ws.Cells.Item(32, 13).Value = -33434.695  # M32: -34476.406 -> -33434.695
# Row 122
$ws.Cells.Item(122, 8).Value = 1676.5834  # H122: 1690.68 -> 1676.5834
$ws.Cells.Item(122, 9).Value = 1555.5454  # I122: 1583.9048 -> 1555.5454
$ws.Cells.Item(122, 10).Value = 3008  # J122: 2251.25 -> 3008
$ws.Cells.Item(122, 11).Value = 4666.6362  # K122: 4751.7144 -> 4666.6362
$ws.Cells.Item(122, 12).Value = 9024  # L122: 6753.75 -> 9024
$ws.Cells.Item(122, 13).Value = -2216.6362  # M122: -2301.7144 -> -2216.6362
$ws.Cells.Item(122, 14).Value = -13924  # N122: -11653.75 -> -13924
# Row 132
$ws.Cells.Item(132, 8).Value = 1551.9546  # H132: 1445.3733 -> 1551.9546
$ws.Cells.Item(132, 9).Value = 1153  # I132: 1063.1428 -> 1153
$ws.Cells.Item(132, 11).Value = 3459  # K132: 3189.4284 -> 3459
$ws.Cells.Item(132, 13).Value = -929  # M132: -659.4284000000002 -> -929

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 2445.7144  # H16: 2153.158 -> 2445.7144
$ws.Cells.Item(16, 9).Value = 1621  # I16: 1509.8572 -> 1621
$ws.Cells.Item(16, 10).Value = 3545.3333  # J16: 3954.4 -> 3545.3333
$ws.Cells.Item(16, 11).Value = 1621  # K16: 1509.8572 -> 1621
$ws.Cells.Item(16, 12).Value = 3545.3333  # L16: 3954.4 -> 3545.3333
$ws.Cells.Item(16, 13).Value = -1334  # M16: -1222.8572 -> -1334
$ws.Cells.Item(16, 14).Value = -4119.3333  # N16: -4528.4 -> -4119.3333
# Row 108
$ws.Cells.Item(108, 8).Value = 40250  # H108: 42750 -> 40250
$ws.Cells.Item(108, 9).Value = 0  # I108: 40000 -> 0
$ws.Cells.Item(108, 10).Value = 40250  # J108: 45500 -> 40250
$ws.Cells.Item(108, 11).Value = 0  # K108: 40000 -> 0
$ws.Cells.Item(108, 12).Value = 40250  # L108: 45500 -> 40250
$ws.Cells.Item(108, 13).ClearContents()  # M108: -36160 -> (blank)
$ws.Cells.Item(108, 14).Value = -47930  # N108: -53180 -> -47930
# Row 113
$ws.Cells.Item(113, 8).Value = 2445.7144  # H113: 2153.158 -> 2445.7144
$ws.Cells.Item(113, 9).Value = 1621  # I113: 1509.8572 -> 1621
$ws.Cells.Item(113, 10).Value = 3545.3333  # J113: 3954.4 -> 3545.3333
$ws.Cells.Item(113, 11).Value = 1621  # K113: 1509.8572 -> 1621
$ws.Cells.Item(113, 12).Value = 3545.3333  # L113: 3954.4 -> 3545.3333
$ws.Cells.Item(113, 13).Value = 549  # M113: 660.1428000000001 -> 549
$ws.Cells.Item(113, 14).Value = -7885.3333  # N113: -8294.4 -> -7885.3333
# Row 132
$ws.Cells.Item(132, 8).Value = 44339  # H132: 40851.88 -> 44339
$ws.Cells.Item(132, 9).Value = 53500.42  # I132: 48476.57 -> 53500.42
$ws.Cells.Item(132, 11).Value = 160501.26  # K132: 145429.71 -> 160501.26
$ws.Cells.Item(132, 13).Value = -157971.26  # M132: -142899.71 -> -157971.26
# Row 134
$ws.Cells.Item(134, 8).Value = 1552.1163  # H134: 1539.7 -> 1552.1163
$ws.Cells.Item(134, 9).Value = 1235.919  # I134: 1265.3235 -> 1235.919
$ws.Cells.Item(134, 10).Value = 3502  # J134: 3094.5 -> 3502
$ws.Cells.Item(134, 11).Value = 3707.757000000001  # K134: 3795.9705 -> 3707.757000000001
$ws.Cells.Item(134, 12).Value = 10506  # L134: 9283.5 -> 10506
$ws.Cells.Item(134, 13).Value = -1172.757000000001  # M134: -1260.9705 -> -1172.757000000001
$ws.Cells.Item(134, 14).Value = -15576  # N134: -14353.5 -> -15576

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Cells.Item(9, 8).Value = 333333540  # H9: 500000300 -> 333333540
$ws.Cells.Item(9, 9).Value = 300  # I9: 0 -> 300
$ws.Cells.Item(9, 10).Value = 500000130  # J9: 500000300 -> 500000130
$ws.Cells.Item(9, 11).Value = 900  # K9: 0 -> 900
$ws.Cells.Item(9, 12).Value = 1500000390  # L9: 1500000900 -> 1500000390
$ws.Cells.Item(9, 13).Value = -676  # M9: None -> -676
$ws.Cells.Item(9, 14).Value = -1500000838  # N9: -1500001348 -> -1500000838
# Row 69
$ws.Cells.Item(69, 8).Value = 4241.5835  # H69: 4620.8335 -> 4241.5835
$ws.Cells.Item(69, 9).Value = 449  # I69: 0 -> 449
$ws.Cells.Item(69, 10).Value = 4586.364  # J69: 4620.8335 -> 4586.364
$ws.Cells.Item(69, 11).Value = 1347  # K69: 0 -> 1347
$ws.Cells.Item(69, 12).Value = 13759.092  # L69: 13862.5005 -> 13759.092
$ws.Cells.Item(69, 13).Value = -536  # M69: None -> -536
$ws.Cells.Item(69, 14).Value = -15381.092  # N69: -15484.5005 -> -15381.092
# Row 72
$ws.Cells.Item(72, 8).Value = 4241.5835  # H72: 4620.8335 -> 4241.5835
$ws.Cells.Item(72, 9).Value = 449  # I72: 0 -> 449
$ws.Cells.Item(72, 10).Value = 4586.364  # J72: 4620.8335 -> 4586.364
$ws.Cells.Item(72, 11).Value = 4041  # K72: 0 -> 4041
$ws.Cells.Item(72, 12).Value = 41277.276  # L72: 41587.5015 -> 41277.276
$ws.Cells.Item(72, 13).Value = 15  # M72: None -> 15
$ws.Cells.Item(72, 14).Value = -49389.276  # N72: -49699.5015 -> -49389.276
# Row 137
$ws.Cells.Item(137, 8).Value = 3787.7778  # H137: 3614.7368 -> 3787.7778
$ws.Cells.Item(137, 9).Value = 3746.7144  # I137: 3340.875 -> 3746.7144
$ws.Cells.Item(137, 11).Value = 11240.1432  # K137: 10022.625 -> 11240.1432
$ws.Cells.Item(137, 13).Value = -6140.143199999999  # M137: -4922.625 -> -6140.143199999999

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 98
$ws.Cells.Item(98, 8).Value = 22998.5  # H98: 23110 -> 22998.5
$ws.Cells.Item(98, 10).Value = 22998.5  # J98: 23110 -> 22998.5
$ws.Cells.Item(98, 12).Value = 22998.5  # L98: 23110 -> 22998.5
$ws.Cells.Item(98, 14).Value = -28988.5  # N98: -29100 -> -28988.5
# Row 132
$ws.Cells.Item(132, 8).Value = 2393  # H132: 2031.0883 -> 2393
$ws.Cells.Item(132, 9).Value = 1723.7059  # I132: 1482.1154 -> 1723.7059
$ws.Cells.Item(132, 11).Value = 5171.1177  # K132: 4446.3462 -> 5171.1177
$ws.Cells.Item(132, 13).Value = -2641.1177  # M132: -1916.3462 -> -2641.1177
# Row 134
$ws.Cells.Item(134, 8).Value = 44000  # H134: 45666.332 -> 44000
$ws.Cells.Item(134, 9).Value = 0  # I134: 48999 -> 0
$ws.Cells.Item(134, 11).Value = 0  # K134: 146997 -> 0
$ws.Cells.Item(134, 13).ClearContents()  # M134: -144462 -> (blank)

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 3134.7827  # H46: 3035.4167 -> 3134.7827
$ws.Cells.Item(46, 9).Value = 820  # I46: 802.5 -> 820
$ws.Cells.Item(46, 11).Value = 820  # K46: 802.5 -> 820
$ws.Cells.Item(46, 13).Value = -632  # M46: -614.5 -> -632
# Row 122
$ws.Cells.Item(122, 8).Value = 5461.4165  # H122: 3136.205 -> 5461.4165
$ws.Cells.Item(122, 9).Value = 3814  # I122: 2486.258 -> 3814
$ws.Cells.Item(122, 10).Value = 7767.8  # J122: 5654.75 -> 7767.8
$ws.Cells.Item(122, 11).Value = 11442  # K122: 7458.773999999999 -> 11442
$ws.Cells.Item(122, 12).Value = 23303.4  # L122: 16964.25 -> 23303.4
$ws.Cells.Item(122, 13).Value = -8992  # M122: -5008.773999999999 -> -8992
$ws.Cells.Item(122, 14).Value = -28203.4  # N122: -21864.25 -> -28203.4
# Row 132
$ws.Cells.Item(132, 8).Value = 3610.6667  # H132: 2685.6155 -> 3610.6667
$ws.Cells.Item(132, 9).Value = 0  # I132: 1082.4 -> 0
$ws.Cells.Item(132, 10).Value = 3610.6667  # J132: 3687.625 -> 3610.6667
$ws.Cells.Item(132, 11).Value = 0  # K132: 3247.2 -> 0
$ws.Cells.Item(132, 12).Value = 10832.0001  # L132: 11062.875 -> 10832.0001
$ws.Cells.Item(132, 13).ClearContents()  # M132: -717.2000000000003 -> (blank)
$ws.Cells.Item(132, 14).Value = -15892.0001  # N132: -16122.875 -> -15892.0001

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 1641.0358  # H132: 1676.5186 -> 1641.0358
$ws.Cells.Item(132, 9).Value = 1188.1428  # I132: 1213.4 -> 1188.1428
$ws.Cells.Item(132, 11).Value = 3564.4284  # K132: 3640.2 -> 3564.4284
$ws.Cells.Item(132, 13).Value = -1034.4284  # M132: -1110.2 -> -1034.4284
# Row 136
$ws.Cells.Item(136, 8).Value = 17542.666  # H136: 18557.412 -> 17542.666
$ws.Cells.Item(136, 9).Value = 18576.879  # I136: 19756.549 -> 18576.879
$ws.Cells.Item(136, 11).Value = 55730.637  # K136: 59269.647 -> 55730.637
$ws.Cells.Item(136, 13).Value = -53180.637  # M136: -56719.647 -> -53180.637

Write-Host "Applied scheduled data refresh updates to all sheets."
